$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1904761904761905
$ws.Range("C2").Value = 0.5238095238095238
$ws.Range("J2").Value = 0.04761904761904762
$ws.Range("P2").Value = 0.1904761904761905
$ws.Range("S2").Value = 0.04761904761904762
$ws.Range("B3").Value = 0.08333333333333333
$ws.Range("J3").Value = 0.08333333333333333
$ws.Range("P3").Value = 0.6666666666666666
$ws.Range("S3").Value = 0.1666666666666667
$ws.Range("J4").Value = 0.1111111111111111
$ws.Range("P4").Value = 0.6666666666666666
$ws.Range("S4").Value = 0.2222222222222222
$ws.Range("B6").Value = 0.1428571428571428
$ws.Range("D6").Value = 0.04761904761904762
$ws.Range("F6").Value = 0.04761904761904762
$ws.Range("J6").Value = 0.1428571428571428
$ws.Range("Q6").Value = 0.1904761904761905
$ws.Range("R6").Value = 0.09523809523809523
$ws.Range("S6").Value = 0.3333333333333333
$ws.Range("D7").Value = 0.2222222222222222
$ws.Range("F7").Value = 0.1111111111111111
$ws.Range("J7").Value = 0.2222222222222222
$ws.Range("O7").Value = 0.1111111111111111
$ws.Range("Q7").Value = 0.1111111111111111
$ws.Range("R7").Value = 0.1111111111111111
$ws.Range("S7").Value = 0.1111111111111111
$ws.Range("D8").Value = 0.05263157894736842
$ws.Range("J8").Value = 0.1052631578947368
$ws.Range("Q8").Value = 0.1578947368421053
$ws.Range("R8").Value = 0.1052631578947368
$ws.Range("S8").Value = 0.5789473684210527
$ws.Range("D9").Value = 0.1666666666666667
$ws.Range("F9").Value = 0.1666666666666667
$ws.Range("J9").Value = 0.1666666666666667
$ws.Range("R9").Value = 0.3333333333333333
$ws.Range("S9").Value = 0.1666666666666667
$ws.Range("B10").Value = 0.1101694915254237
$ws.Range("D10").Value = 0.03389830508474576
$ws.Range("E10").Value = 0.00847457627118644
$ws.Range("F10").Value = 0.1016949152542373
$ws.Range("J10").Value = 0.1186440677966102
$ws.Range("O10").Value = 0.00847457627118644
$ws.Range("Q10").Value = 0.2457627118644068
$ws.Range("R10").Value = 0.07627118644067797
$ws.Range("S10").Value = 0.2966101694915254
$ws.Range("G11").Value = 0.25
$ws.Range("J11").Value = 0.08333333333333333
$ws.Range("K11").Value = 0.25
$ws.Range("L11").Value = 0.4166666666666667
$ws.Range("G12").Value = 0.8
$ws.Range("J12").Value = 0.2
$ws.Range("J15").Value = 0.5882352941176471
$ws.Range("K15").Value = 0.1176470588235294
$ws.Range("O15").Value = 0.1176470588235294
$ws.Range("S15").Value = 0.1764705882352941
$ws.Range("J16").Value = 0.7222222222222222
$ws.Range("K16").Value = 0.05555555555555555
$ws.Range("O16").Value = 0.1666666666666667
$ws.Range("S16").Value = 0.05555555555555555
$ws.Range("F17").Value = 0.02857142857142857
$ws.Range("H17").Value = 0.05714285714285714
$ws.Range("I17").Value = 0.08571428571428572
$ws.Range("J17").Value = 0.5428571428571428
$ws.Range("M17").Value = 0.02857142857142857
$ws.Range("O17").Value = 0.05714285714285714
$ws.Range("J18").Value = 0.75
$ws.Range("K18").Value = 0.0625
$ws.Range("O18").Value = 0.0625
$ws.Range("S18").Value = 0.125
$ws.Range("F19").Value = 0.02564102564102564
$ws.Range("H19").Value = 0.217948717948718
$ws.Range("I19").Value = 0.03846153846153846
$ws.Range("J19").Value = 0.4871794871794872
$ws.Range("K19").Value = 0.0641025641025641
$ws.Range("M19").Value = 0.01282051282051282
$ws.Range("O19").Value = 0.07692307692307693
$ws.Range("S19").Value = 0.07692307692307693
